# Append three paragraphs after the final "*master" paragraph:
#   1) a blank ListParagraph-styled separator paragraph
#   2) a numbered ListParagraph item: "Push code to master" (continues the
#      existing numbered list, numId=1, so it becomes item "8.")
#   3) a ListParagraph (no numbering) command line:
#      ">>git push -u origin master"

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Create three new empty paragraphs at the end of the document, each
# inheriting the style of the paragraph before it (List Paragraph).
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)            # wdCollapseEnd
$r.InsertParagraphAfter()

$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$total = $d.Paragraphs.Count
$pBlank = $d.Paragraphs.Item($total - 2)
$pPush  = $d.Paragraphs.Item($total - 1)
$pCmd   = $d.Paragraphs.Item($total)

# 1) Blank separator paragraph - ListParagraph style, no numbering, no run.
$xmlBlank = "<w:p $wNs><w:pPr><w:pStyle w:val=""ListParagraph""/></w:pPr></w:p>"
$pBlank.Range.InsertXML($xmlBlank)

# 2) Numbered list item "Push code to master" - continues numId 1.
$xmlPush = "<w:p $wNs><w:pPr><w:pStyle w:val=""ListParagraph""/><w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""1""/></w:numPr></w:pPr><w:r><w:t>Push code to master</w:t></w:r></w:p>"
$pPush.Range.InsertXML($xmlPush)

# 3) Command line paragraph ">>git push -u origin master".
$xmlCmd = "<w:p $wNs><w:pPr><w:pStyle w:val=""ListParagraph""/></w:pPr><w:r><w:t>&gt;&gt;git push -u origin master</w:t></w:r></w:p>"
$pCmd.Range.InsertXML($xmlCmd)
